$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "30.169.89"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -0.52%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.864.67"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -0.39%  "
$ws.Range("E4").Value = "  -0.03%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "234.40"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.72%  "
$ws.Range("E6").Value = "  +0.00%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.4680"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.67%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.2867"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -0.74%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.06497"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -2.09%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "21.18"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -2.56%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.07755"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -3.61%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "1.872.80"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +0.02%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "94.01"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -3.57%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.6848"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -0.58%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "5.057"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -1.78%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "269.52"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -1.05%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "30.159.29"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -0.53%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "13.34"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -6.05%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.000007652"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -0.04%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "2.110.64"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -0.31%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.02%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "5.163"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -3.11%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "6.109"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -1.84%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "9.350"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -0.30%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "165.80"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -1.45%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "18.57"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -2.34%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "1.896"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -3.31%  "
$ws.Range("E29").Value = "  -0.78%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "0.09930"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -0.43%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "1.451"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -0.99%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "4.235"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -3.30%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "4.016"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -1.91%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.04689"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -0.59%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.120"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -1.91%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.6899"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -1.94%  "
$ws.Range("E37").Value = "  -0.40%  "
$ws.Range("E38").Value = "  -2.74%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "2.760"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +4.03%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "6.345"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +0.29%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "71.46"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -1.88%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +0.01%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "1.898"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -3.42%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.8346"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -1.18%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "102.39"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -0.80%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.4067"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -2.71%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "937.89"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +0.27%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "9.119"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -2.26%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "6.973"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -2.00%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "34.09"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -1.33%  "
$ws.Range("E51").Value = "  -1.82%  "
